$d = $word.ActiveDocument

$d.Content.Find.Execute("90-1=", $true, $false, $false, $false, $false, $true, 1, $false, "30+36=", 2) | Out-Null
$d.Content.Find.Execute("52-48=", $true, $false, $false, $false, $false, $true, 1, $false, "76+19=", 2) | Out-Null
$d.Content.Find.Execute("90-16=", $true, $false, $false, $false, $false, $true, 1, $false, "46+13=", 2) | Out-Null
$d.Content.Find.Execute("15-2=", $true, $false, $false, $false, $false, $true, 1, $false, "51-4=", 2) | Out-Null
$d.Content.Find.Execute("27+53=", $true, $false, $false, $false, $false, $true, 1, $false, "17+7=", 2) | Out-Null
$d.Content.Find.Execute("30+3=", $true, $false, $false, $false, $false, $true, 1, $false, "51-1=", 2) | Out-Null
$d.Content.Find.Execute("78-43=", $true, $false, $false, $false, $false, $true, 1, $false, "99-24=", 2) | Out-Null
$d.Content.Find.Execute("71-42=", $true, $false, $false, $false, $false, $true, 1, $false, "17+49=", 2) | Out-Null
$d.Content.Find.Execute("63-18=", $true, $false, $false, $false, $false, $true, 1, $false, "48-17=", 2) | Out-Null
$d.Content.Find.Execute("89-65=", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=", 2) | Out-Null
$d.Content.Find.Execute("61+7=", $true, $false, $false, $false, $false, $true, 1, $false, "85-16=", 2) | Out-Null
$d.Content.Find.Execute("88-37=", $true, $false, $false, $false, $false, $true, 1, $false, "52+10=", 2) | Out-Null
$d.Content.Find.Execute("60+6=", $true, $false, $false, $false, $false, $true, 1, $false, "51+6=", 2) | Out-Null
$d.Content.Find.Execute("1+97=", $true, $false, $false, $false, $false, $true, 1, $false, "8+24=", 2) | Out-Null
$d.Content.Find.Execute("19+13=", $true, $false, $false, $false, $false, $true, 1, $false, "51+28=", 2) | Out-Null
$d.Content.Find.Execute("56-29=", $true, $false, $false, $false, $false, $true, 1, $false, "64-26=", 2) | Out-Null
$d.Content.Find.Execute("47-11=", $true, $false, $false, $false, $false, $true, 1, $false, "3+30=", 2) | Out-Null
$d.Content.Find.Execute("4+33=", $true, $false, $false, $false, $false, $true, 1, $false, "69-54=", 2) | Out-Null
$d.Content.Find.Execute("18+58=", $true, $false, $false, $false, $false, $true, 1, $false, "5+85=", 2) | Out-Null
$d.Content.Find.Execute("89-27=", $true, $false, $false, $false, $false, $true, 1, $false, "61-1=", 2) | Out-Null
$d.Content.Find.Execute("41-8=", $true, $false, $false, $false, $false, $true, 1, $false, "21-21=", 2) | Out-Null
$d.Content.Find.Execute("40-25=", $true, $false, $false, $false, $false, $true, 1, $false, "64-49=", 2) | Out-Null
$d.Content.Find.Execute("33+9=", $true, $false, $false, $false, $false, $true, 1, $false, "61+15=", 2) | Out-Null
$d.Content.Find.Execute("99-29=", $true, $false, $false, $false, $false, $true, 1, $false, "83-11=", 2) | Out-Null
$d.Content.Find.Execute("58+30=", $true, $false, $false, $false, $false, $true, 1, $false, "15+70=", 2) | Out-Null
$d.Content.Find.Execute("99-70=", $true, $false, $false, $false, $false, $true, 1, $false, "73-50=", 2) | Out-Null
$d.Content.Find.Execute("57+36=", $true, $false, $false, $false, $false, $true, 1, $false, "55+31=", 2) | Out-Null
$d.Content.Find.Execute("56-37=", $true, $false, $false, $false, $false, $true, 1, $false, "82-35=", 2) | Out-Null
$d.Content.Find.Execute("4+83=", $true, $false, $false, $false, $false, $true, 1, $false, "67-30=", 2) | Out-Null
$d.Content.Find.Execute("73+21=", $true, $false, $false, $false, $false, $true, 1, $false, "39-25=", 2) | Out-Null
$d.Content.Find.Execute("68+21=", $true, $false, $false, $false, $false, $true, 1, $false, "36+37=", 2) | Out-Null
$d.Content.Find.Execute("67+2=", $true, $false, $false, $false, $false, $true, 1, $false, "54-54=", 2) | Out-Null
$d.Content.Find.Execute("72-53=", $true, $false, $false, $false, $false, $true, 1, $false, "72-65=", 2) | Out-Null
$d.Content.Find.Execute("46+17=", $true, $false, $false, $false, $false, $true, 1, $false, "27+0=", 2) | Out-Null
$d.Content.Find.Execute("75-65=", $true, $false, $false, $false, $false, $true, 1, $false, "42-34=", 2) | Out-Null
$d.Content.Find.Execute("92-43=", $true, $false, $false, $false, $false, $true, 1, $false, "41+28=", 2) | Out-Null
$d.Content.Find.Execute("21+58=", $true, $false, $false, $false, $false, $true, 1, $false, "79-39=", 2) | Out-Null
$d.Content.Find.Execute("90-76=", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=", 2) | Out-Null
$d.Content.Find.Execute("96-27=", $true, $false, $false, $false, $false, $true, 1, $false, "25+73=", 2) | Out-Null
$d.Content.Find.Execute("8+37=", $true, $false, $false, $false, $false, $true, 1, $false, "46+44=", 2) | Out-Null
$d.Content.Find.Execute("37+57=", $true, $false, $false, $false, $false, $true, 1, $false, "3+39=", 2) | Out-Null
$d.Content.Find.Execute("81+8=", $true, $false, $false, $false, $false, $true, 1, $false, "57+15=", 2) | Out-Null
$d.Content.Find.Execute("13+17=", $true, $false, $false, $false, $false, $true, 1, $false, "42-14=", 2) | Out-Null
$d.Content.Find.Execute("96-68=", $true, $false, $false, $false, $false, $true, 1, $false, "97-44=", 2) | Out-Null
$d.Content.Find.Execute("24+42=", $true, $false, $false, $false, $false, $true, 1, $false, "6+15=", 2) | Out-Null
$d.Content.Find.Execute("60+36=", $true, $false, $false, $false, $false, $true, 1, $false, "36-29=", 2) | Out-Null
$d.Content.Find.Execute("91-40=", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=", 2) | Out-Null
$d.Content.Find.Execute("55+0=", $true, $false, $false, $false, $false, $true, 1, $false, "6+38=", 2) | Out-Null
$d.Content.Find.Execute("17+17=", $true, $false, $false, $false, $false, $true, 1, $false, "33+16=", 2) | Out-Null
$d.Content.Find.Execute("30+2=", $true, $false, $false, $false, $false, $true, 1, $false, "63-7=", 2) | Out-Null
$d.Content.Find.Execute("90-88=", $true, $false, $false, $false, $false, $true, 1, $false, "61-61=", 2) | Out-Null
$d.Content.Find.Execute("8+43=", $true, $false, $false, $false, $false, $true, 1, $false, "40+25=", 2) | Out-Null
$d.Content.Find.Execute("80-25=", $true, $false, $false, $false, $false, $true, 1, $false, "18+7=", 2) | Out-Null
$d.Content.Find.Execute("83-13=", $true, $false, $false, $false, $false, $true, 1, $false, "74+9=", 2) | Out-Null
$d.Content.Find.Execute("16-3=", $true, $false, $false, $false, $false, $true, 1, $false, "36+34=", 2) | Out-Null
$d.Content.Find.Execute("89-12=", $true, $false, $false, $false, $false, $true, 1, $false, "92-25=", 2) | Out-Null
$d.Content.Find.Execute("98+1=", $true, $false, $false, $false, $false, $true, 1, $false, "10+65=", 2) | Out-Null
$d.Content.Find.Execute("38+15=", $true, $false, $false, $false, $false, $true, 1, $false, "18+81=", 2) | Out-Null
$d.Content.Find.Execute("79-4=", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=", 2) | Out-Null
$d.Content.Find.Execute("70+15=", $true, $false, $false, $false, $false, $true, 1, $false, "93-75=", 2) | Out-Null
$d.Content.Find.Execute("78+19=", $true, $false, $false, $false, $false, $true, 1, $false, "1-1=", 2) | Out-Null
$d.Content.Find.Execute("33-31=", $true, $false, $false, $false, $false, $true, 1, $false, "55+17=", 2) | Out-Null
$d.Content.Find.Execute("61-27=", $true, $false, $false, $false, $false, $true, 1, $false, "33+25=", 2) | Out-Null
$d.Content.Find.Execute("74-34=", $true, $false, $false, $false, $false, $true, 1, $false, "63-56=", 2) | Out-Null
$d.Content.Find.Execute("1+93=", $true, $false, $false, $false, $false, $true, 1, $false, "85-23=", 2) | Out-Null
$d.Content.Find.Execute("45-36=", $true, $false, $false, $false, $false, $true, 1, $false, "92-37=", 2) | Out-Null
$d.Content.Find.Execute("50+15=", $true, $false, $false, $false, $false, $true, 1, $false, "63-60=", 2) | Out-Null
$d.Content.Find.Execute("74+17=", $true, $false, $false, $false, $false, $true, 1, $false, "98-42=", 2) | Out-Null
$d.Content.Find.Execute("21-20=", $true, $false, $false, $false, $false, $true, 1, $false, "35-14=", 2) | Out-Null
$d.Content.Find.Execute("48+16=", $true, $false, $false, $false, $false, $true, 1, $false, "51+11=", 2) | Out-Null
$d.Content.Find.Execute("85-25=", $true, $false, $false, $false, $false, $true, 1, $false, "10+18=", 2) | Out-Null
$d.Content.Find.Execute("50+2=", $true, $false, $false, $false, $false, $true, 1, $false, "83-36=", 2) | Out-Null
$d.Content.Find.Execute("3+23=", $true, $false, $false, $false, $false, $true, 1, $false, "85-12=", 2) | Out-Null
$d.Content.Find.Execute("6+51=", $true, $false, $false, $false, $false, $true, 1, $false, "79-54=", 2) | Out-Null
$d.Content.Find.Execute("39-7=", $true, $false, $false, $false, $false, $true, 1, $false, "68-41=", 2) | Out-Null
$d.Content.Find.Execute("90-37=", $true, $false, $false, $false, $false, $true, 1, $false, "58-33=", 2) | Out-Null
$d.Content.Find.Execute("89-1=", $true, $false, $false, $false, $false, $true, 1, $false, "18+46=", 2) | Out-Null
$d.Content.Find.Execute("75+5=", $true, $false, $false, $false, $false, $true, 1, $false, "28+26=", 2) | Out-Null
$d.Content.Find.Execute("77-23=", $true, $false, $false, $false, $false, $true, 1, $false, "26+57=", 2) | Out-Null
$d.Content.Find.Execute("55+25=", $true, $false, $false, $false, $false, $true, 1, $false, "96-34=", 2) | Out-Null
$d.Content.Find.Execute("83-71=", $true, $false, $false, $false, $false, $true, 1, $false, "99-83=", 2) | Out-Null
$d.Content.Find.Execute("97-75=", $true, $false, $false, $false, $false, $true, 1, $false, "85-4=", 2) | Out-Null
$d.Content.Find.Execute("35+49=", $true, $false, $false, $false, $false, $true, 1, $false, "27+44=", 2) | Out-Null
$d.Content.Find.Execute("23+63=", $true, $false, $false, $false, $false, $true, 1, $false, "25+17=", 2) | Out-Null
$d.Content.Find.Execute("9+57=", $true, $false, $false, $false, $false, $true, 1, $false, "71+3=", 2) | Out-Null
$d.Content.Find.Execute("1+92=", $true, $false, $false, $false, $false, $true, 1, $false, "97-8=", 2) | Out-Null
$d.Content.Find.Execute("98-24=", $true, $false, $false, $false, $false, $true, 1, $false, "68-20=", 2) | Out-Null
$d.Content.Find.Execute("67-57=", $true, $false, $false, $false, $false, $true, 1, $false, "50+40=", 2) | Out-Null
$d.Content.Find.Execute("3+11=", $true, $false, $false, $false, $false, $true, 1, $false, "4+69=", 2) | Out-Null
$d.Content.Find.Execute("10+81=", $true, $false, $false, $false, $false, $true, 1, $false, "35+52=", 2) | Out-Null
$d.Content.Find.Execute("24-23=", $true, $false, $false, $false, $false, $true, 1, $false, "86-25=", 2) | Out-Null
$d.Content.Find.Execute("11+57=", $true, $false, $false, $false, $false, $true, 1, $false, "64+10=", 2) | Out-Null
$d.Content.Find.Execute("53-1=", $true, $false, $false, $false, $false, $true, 1, $false, "69-5=", 2) | Out-Null
$d.Content.Find.Execute("64-1=", $true, $false, $false, $false, $false, $true, 1, $false, "55-33=", 2) | Out-Null
$d.Content.Find.Execute("9+41=", $true, $false, $false, $false, $false, $true, 1, $false, "78+5=", 2) | Out-Null
$d.Content.Find.Execute("19-5=", $true, $false, $false, $false, $false, $true, 1, $false, "60-5=", 2) | Out-Null
$d.Content.Find.Execute("15+52=", $true, $false, $false, $false, $false, $true, 1, $false, "66-7=", 2) | Out-Null
$d.Content.Find.Execute("44-25=", $true, $false, $false, $false, $false, $true, 1, $false, "21+24=", 2) | Out-Null
$d.Content.Find.Execute("81-27=", $true, $false, $false, $false, $false, $true, 1, $false, "41+48=", 2) | Out-Null
$d.Content.Find.Execute("85-22=", $true, $false, $false, $false, $false, $true, 1, $false, "58-31=", 2) | Out-Null
